$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the <PHOTO 1> / <PHOTO 2> / <PHOTO 3> placeholder text from row 10
$ws.Range("A10:C10").ClearContents()

# Row 10 grows taller now that the photo placeholders are gone
$ws.Rows.Item(10).RowHeight = 133.8

# Update the view: zoom to 85% in Page Break Preview and move the selection to A9:C9
$excel.ActiveWindow.Zoom = 85
$ws.Range("A9:C9").Select()
